# edit.ps1 - applies the LOT2037.docx content revision (2025-05-09 build) via Word COM interop
$d = $word.ActiveDocument

# 1) Update 'Ativacao' (activation) date from 2018 to 2025
$ok1 = $d.Content.Find.Execute('Ativação: 01/01/2018', $true, $false, $false, $false, $false, $true, 1, $false, 'Ativação: 01/01/2025', 2)
if (-not $ok1) { Write-Host "WARNING: step ok1 did not find its target text" }

# 2) Replace the Portuguese 'Objetivos' paragraph with the new objective text
$ok2 = $d.Content.Find.Execute('Apresentar as modernas técnicas de instrumentação, monitoramento e controle de bioprocessos em geral, bem como estudar os fundamentos das diversas técnicas apresentadas e o princípio de funcionamento dos diferentes equipamentos utilizados nas medidas.', $true, $false, $false, $false, $false, $true, 1, $false, 'Desenvolver nos discentes as competências e habilidades necessárias para aplicar conhecimentos científicos, tecnológicos e de engenharia na concepção, projeto, instalação, otimização, supervisão e avaliação crítica da operação de bioprocessos, com ênfase em: 1) Abordar as categorias de instrumentos e o uso de símbolos ISA em diagramas P&ID. 2) Explorar a medição de pressão, nível, vazão e temperatura, bem como capacitar para o uso de instrumentos na monitoração e controle de parâmetros em biorreatores, incorporando diversas estratégias de controle.', 2)
if (-not $ok2) { Write-Host "WARNING: step ok2 did not find its target text" }

# 3) The English (italic) objectives paragraph exists but is empty in the source document;
#    add its translated text run (keeps the existing italic formatting)
$obj_en_inserted = $false
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text.TrimEnd([char]13, [char]7) -eq "" -and $r.Font.Italic -eq -1) {
        $r.Text = 'Develop in students the competencies and skills necessary to apply scientific, technological, and engineering knowledge in the design, project, installation, optimization, supervision, and critical evaluation of bioprocess operations, with an emphasis on: 1) Addressing the categories of instruments and the use of ISA symbols in P&ID diagrams. 2) Exploring the measurement of pressure, level, flow, and temperature, as well as equipping students to use instruments for monitoring and controlling parameters in bioreactors, incorporating various control strategies.'
        $obj_en_inserted = $true
        break
    }
}
if (-not $obj_en_inserted) { Write-Host "WARNING: empty italic objectives run not found" }

# 4) Replace the Portuguese 'Programa resumido' paragraph
$ok4 = $d.Content.Find.Execute('Introdução. Conceitos básicos de instrumentação para controle de processos. Instrumentos de medição de pressão. Dispositivos de medição de nível. Instrumentos de medição de vazão. Dispositivos de medição de temperatura. Sensores online para medição de meios, células e gases. Introdução à teoria de controle. Controles aplicados a bioprocessos.', $true, $false, $false, $false, $false, $true, 1, $false, 'Importância da instrumentação e controle em bioprocessos industriais. Categorias de instrumentos e uso de símbolos ISA em diagramas P&ID. Instrumentos de medição de pressão. Dispositivos para medição direta e indireta de nível. Estudo de instrumentos de medição de vazão. Dispositivos para medição de temperatura. Uso de instrumentos para monitorar e controlar parâmetros em biorreatores, com diversas estratégias de controle.', 2)
if (-not $ok4) { Write-Host "WARNING: step ok4 did not find its target text" }

# 5) Replace the English 'Programa resumido' (italic) paragraph
$ok5 = $d.Content.Find.Execute('Introduction. Basic concepts of instrumentation for process control. Pressure measuring instruments. Level measuring devices. Flow measuring instruments. Temperature measuring devices. Online sensors for measurement of media, cells and gases. Introduction to control theory. Controls applied to bioprocesses', $true, $false, $false, $false, $false, $true, 1, $false, 'Importance of instrumentation and control in industrial bioprocesses. Categories of instruments and the use of ISA symbols in P&ID diagrams. Pressure measurement instruments. Devices for direct and indirect level measurement. Study of flow measurement instruments. Devices for temperature measurement. Use of instruments to monitor and control parameters in bioreactors, incorporating various control strategies.', 2)
if (-not $ok5) { Write-Host "WARNING: step ok5 did not find its target text" }

# 6) Replace the Portuguese 'Programa' (detailed syllabus) paragraph
$ok6 = $d.Content.Find.Execute('1.Introdução: aspectos gerais relativos à instrumentação e controle de bioprocessos.2.Conceitos básicos de instrumentação para controle de processos: definições dos elementos em uma malha de controle. Características gerais de instrumentos: classes de instrumentos e definições. Identificação e símbolos de instrumentos: padronização ISA, exemplos de simbologia.3.Instrumentos de medição de pressão: manômetros, diafragmas, cápsulas e foles, tubos de Bourdon e outros sensores de pressão.4.Dispositivos de medição de nível: medição direta de nível e medição indireta de nível.5.Instrumentos de medição de vazão: medidores deprimogênios, medidores lineares, medidores volumétricos e outros.6.Dispositivos de medição de temperatura: termômetros, termômetros com mola de pressão, dispositivos de temperatura de resistência, termistores, termopares e outros.7.Instrumentos de medição de pH, potencial redox, pressões parciais de oxigênio dissolvido e gás carbônico. Medição de potência de agitação e velocidade do impelidor. Sensores online para propriedades celulares e determinação da concentração total de biomassa. Determinação da concentração de biomassa ativa ou viável.8. Introdução à teoria de controle: principais problemas para o controle de bioprocessos. Definições básicas (controle manual, controle por realimentação - feedback, controle por antecipação - feedforward, ganho e atraso), componentes de um sistema de controle (sensor/transmissor e controlador e elementos finais de controle), ações de controle Liga-desliga (on-off), auto-operado, proporcional (P), proporcional-integral (PI), proporcional-derivativa (PD), proporcional-integral-derivativa (PDI). Introdução à interface de comunicação.9.Controles aplicados a bioprocessos: controle em malha aberta, controle por sistema regulatório, controle em cascata, controle por pré-alimentação e controle seguidor de trajetória e outros (controle ótimo, sistema adaptativo e sistema de controle por aprendizado).', $true, $false, $false, $false, $false, $true, 1, $false, '1) Introdução: apresentar a visão geral do uso de instrumentação e controle em bioprocessos, destacando a importância desses instrumentos no setor industrial.2) Conceitos Básicos de Instrumentação: explorar as categorias de instrumentos em malhas de controle e utilizar os símbolos padrão ISA para diagramas P&ID.3) Medição de Pressão: analisar os instrumentos utilizados para medir pressão, como manômetros, diafragmas, cápsulas e foles, tubos de Bourdon e outros tipos de sensores de pressão.4) Medição de Nível: estudar dispositivos para medição de nível, incluindo métodos de medição direta, como réguas e visores de nível, e técnicas de medição indireta, como transdutores de nível, sensores capacitivos, de radar e ultrassônicos.5) Medição de Vazão: estudar instrumentos utilizados para medir vazão, englobando medidores deprimogênios, medidores lineares, medidores volumétricos e outras tecnologias relevantes para a medição de vazão.6) Medição de Temperatura: abordar os dispositivos para medição de temperatura, incluindo termômetros comuns, termômetros com mola de pressão, dispositivos de temperatura de resistência (RTDs), termistores, termopares e outros tipos de sensores de temperatura.7) Monitoramento e Controle de Biorreatores: focar no uso de instrumentos para medir e ajustar parâmetros em biorreatores, como pH, oxigênio dissolvido, gás carbônico e outros. Discutir diferentes métodos de controle, desde técnicas manuais até avançadas, incluindo controle em cascata, estratégias PID e sistemas de controle adaptativo, para aprimorar a eficiência dos bioprocessos.', 2)
if (-not $ok6) { Write-Host "WARNING: step ok6 did not find its target text" }

# 7) Replace the English 'Programa' (detailed syllabus, italic) paragraph
$ok7 = $d.Content.Find.Execute('1.Introduction: general aspects related to the instrumentation and control of bioprocesses.2.Basic concepts of instrumentation for process control: definitions of the elements in a control mesh. General instrument characteristics: instrument classes and definitions. Instrument identification and symbols: ISA standardization, symbology examples.3.Pressure measuring instruments: manometers, diaphragms, capsule and bellows, Bourdon tubes and other pressure sensors.4.Level measuring devices: direct level measurement and indirect level measurement.5.Flow measurement instruments: pressure meters, linear meters, volumetric meters and others.6.Temperature measuring devices: thermometers, pressure spring thermometers, temperature resistance devices, thermistors, thermocouples and others.7.pH measuring instruments, redox potential, partial pressures of dissolved oxygen and carbon dioxide. Measurement of agitation power and impeller speed. Online sensors for cell properties and determination of total biomass concentration. Determination of active or viable biomass concentration.8.Introduction to control theory: main problems for the control of bioprocesses. Basic settings (manual control, feedback control, feedforward control, gain and delay), components of a control system (sensor / transmitter and controller and final control elements), control actions On/Off (P), proportional-integral (PI), proportional-derivative (PD), proportional-integral-derivative (PDI). Introduction to the communication interface.9.Controls applied to bioprocesses: open loop control, control by regulatory system, cascade control, pre-feed control and trajectory tracking control and others (optimal control, adaptive system and learning control system).', $true, $false, $false, $false, $false, $true, 1, $false, '1)Introduction: Provide an overview of the use of instrumentation and control in bioprocesses, highlighting the importance of these instruments in the industrial sector.2)Basic Concepts of Instrumentation: Explore the categories of instruments in control loops and use the standard ISA symbols for P&ID diagrams.3)Pressure Measurement: Analyze the instruments used to measure pressure, such as manometers, diaphragms, capsules, bellows, Bourdon tubes, and other types of pressure sensors.4)Level Measurement: Study devices for level measurement, including direct measurement methods like gauges and level sight glasses, and indirect measurement techniques such as level transducers, capacitive sensors, radar sensors, and ultrasonic sensors.5)Flow Measurement: Study instruments used to measure flow, including differential pressure meters, linear meters, volumetric meters, and other relevant flow measurement technologies.6)Temperature Measurement: Cover devices for temperature measurement, including common thermometers, pressure spring thermometers, resistance temperature devices (RTDs), thermistors, thermocouples, and other types of temperature sensors.7)Bioreactor Monitoring and Control: Focus on the use of instruments to measure and adjust parameters in bioreactors, such as pH, dissolved oxygen, carbon dioxide, and others. Discuss different control methods, from manual techniques to advanced strategies, including cascade control, PID strategies, and adaptive control systems, to enhance the efficiency of bioprocesses.', 2)
if (-not $ok7) { Write-Host "WARNING: step ok7 did not find its target text" }

# 8) Replace the grading 'Criterio' formula sentence
$ok8 = $d.Content.Find.Execute('Média do período normal = (P1 + P2 + T)/3', $true, $false, $false, $false, $false, $true, 1, $false, 'A média do período normal será calculada pela fórmula: Média do período normal = (P1 + P2 + T)/3.', 2)
if (-not $ok8) { Write-Host "WARNING: step ok8 did not find its target text" }

# 9) Merge the two 'Norma de recuperacao' sentences (originally split across a line break)
#    into a single run with revised wording and no <w:br/> between them
$ok9 = $d.Content.Find.Execute('Aos alunos que obtiverem média igual ou maior que 3,0 e menor que 5,0 será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2.Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0.', $true, $false, $false, $false, $false, $true, 1, $false, 'Para os alunos que alcançarem média igual ou superior a 3,0 e inferior a 5,0, será oferecido um programa de recuperação, avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota da prova final)/2. Serão aprovados os alunos que obtiverem média final igual ou superior a 5,0.', 2)
if (-not $ok9) { Write-Host "WARNING: step ok9 did not find its target text" }

# 10) Replace the Bibliografia body paragraph. The original paragraph is a single run made
#     of many short <w:t> fragments separated by <w:br/> line breaks; Find/Replace with a
#     literal vertical-tab encoded pattern would be unwieldy for ~20 breaks, so instead the
#     whole paragraph range (minus its trailing paragraph mark) is located by its distinctive
#     leading text and its contents are replaced outright, collapsing it into one clean run.
$biblio_replaced = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Bibliografia" + [char]11)) {
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Text = 'BEGA, Egidio Alberto (Editor); IBP (Autor). Instrumentação Industrial. 1. ed. São Paulo: Editora Érica, 2011. 694 p.Coughanowr, D.; LeBlanc, S. Process Systems Analysis and Control. [Capa dura]. Edição Inglês. Editora: McGraw-Hill, 2008. ISBN-10: 007339789X, ISBN-13: 978-0073397894.DORAN, Pauline M. Bioprocess Engineering Principles. [Capa dura]. Edição Inglês. Editora: Academic Press, 1995. ISBN-10: 0122208552, ISBN-13: 978-0122208553.DUNN, William C. Fundamentos de Instrumentação Industrial e Controle de Processos. Porto Alegre: Bookman Editora, 2013. 336 p.FRANCHI, Claiton Moro (Org.). Instrumentação de Processos Industriais. 1. ed. São Paulo: Editora Érica, 2014. 336 p. ISBN 9788536512174.SCHMIDELL, W. et al. Biotecnologia Industrial - Engenharia Bioquímica (Vol 2). São Paulo: Edgard Blucher Ltda, 2001.'
        $biblio_replaced = $true
        break
    }
}
if (-not $biblio_replaced) { Write-Host "WARNING: Bibliografia paragraph not found" }
